{"js": "// The document's single paragraph reads \"Version 2.\" (the \"2\" and the\n// final \".\" are in separate runs, with a \"_GoBack\" bookmark sitting\n// between them). The edit changes the displayed version number so the\n// paragraph reads \"Version 1.\" instead, leaving the \"Version\" text and\n// the bookmark untouched.\nconst body = context.document.body;\n\n// Prefer a tight match on \" 2.\" so only the number + trailing period are\n// touched (mirrors the minimal text change: \"2.\" -> \"1.\").\nlet results = body.search(\" 2.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  // Fallback in case spacing/formatting differs from what we expect.\n  results = body.search(\"Version 2.\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"Version 1.\", Word.InsertLocation.replace);\n  }\n} else {\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\" 1.\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document's single paragraph reads \"Version 2.\" \u2014 the final \".\"\n# sits in its own run just after the \"_GoBack\" bookmark. The edit bumps\n# the version number down to 1 and, in doing so, the trailing \".\"\n# collapses into the run that holds the number (so the bookmark ends up\n# as the very last thing in the paragraph instead of sitting before a\n# standalone \".\" run).\n$d = $word.ActiveDocument\n\n# Step 1: drop the standalone trailing \".\" (the very last character in\n# the story, right before the paragraph mark) so the bookmark becomes\n# the last element in the paragraph once the number is rewritten.\n$periodRange = $d.Range($d.Content.End - 2, $d.Content.End - 1)\nif ($periodRange.Text -eq \".\") {\n    $periodRange.Delete()\n}\n\n# Step 2: change the version number and restore the trailing period,\n# using a Find/Replace that targets only the digit so it doesn't span\n# (and thus doesn't disturb) the bookmark sitting right before it.\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.Text = \"2\"\n$find.Replacement.Text = \"1.\"\n$found = $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n\nif (-not $found) {\n    # Fallback: text/spacing differed from what we expected above, so\n    # fall back to a plain whole-text replace for the full \"Version 2.\"\n    # / \"Version 2\" phrase.\n    $d = $word.ActiveDocument\n    $find2 = $d.Content.Find\n    $find2.Text = \"Version 2.\"\n    $find2.Replacement.Text = \"Version 1.\"\n    $found2 = $find2.Execute([ref]$find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find2.Replacement.Text, 2)\n\n    if (-not $found2) {\n        $d = $word.ActiveDocument\n        $find3 = $d.Content.Find\n        $find3.Text = \"Version 2\"\n        $find3.Replacement.Text = \"Version 1\"\n        $find3.Execute([ref]$find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find3.Replacement.Text, 2)\n    }\n}\n"}
